$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.288.20"
$ws.Range("E2").Value = "  -1.78%  "
$ws.Range("D3").Value = "2.447.76"
$ws.Range("E3").Value = "  -1.34%  "
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").Value = "'568.07"
$ws.Range("E5").Value = "  -1.75%  "
$ws.Range("D6").Value = "'145.63"
$ws.Range("E6").Value = "  -0.99%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("E8").Value = "  -2.43%  "
$ws.Range("E9").Value = "  -1.94%  "
$ws.Range("E10").Value = "  +0.11%  "
$ws.Range("D11").Value = "'5.23"
$ws.Range("E11").Value = "  -1.33%  "
$ws.Range("D12").Value = "'0.349"
$ws.Range("D13").Value = "'28.67"
$ws.Range("E13").Value = "  -1.85%  "
$ws.Range("E14").Value = "  -4.14%  "
$ws.Range("D15").Value = "2.892.03"
$ws.Range("E15").Value = "  -0.91%  "
$ws.Range("D16").Value = "62.316.66"
$ws.Range("E16").Value = "  -1.54%  "
$ws.Range("D17").Value = "2.451.74"
$ws.Range("E17").Value = "  -0.99%  "
$ws.Range("D18").Value = "'7.77"
$ws.Range("E18").Value = "  -2.28%  "
$ws.Range("D19").Value = "'10.73"
$ws.Range("E19").Value = "  -3.50%  "
$ws.Range("D20").Value = "'321.85"
$ws.Range("E20").Value = "  -2.99%  "
$ws.Range("E21").Value = "  -0.26%  "
$ws.Range("E22").Value = "  -3.67%  "
$ws.Range("D23").Value = "'0.999"
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("D24").Value = "'9.81"
$ws.Range("E24").Value = "  +6.28%  "
$ws.Range("E25").Value = "  -2.11%  "
$ws.Range("D26").Value = "'635.66"
$ws.Range("E26").Value = "  -6.21%  "
$ws.Range("D27").Value = "2.587.38"
$ws.Range("E27").Value = "  -0.39%  "
$ws.Range("D28").Value = "0.0₃0954"
$ws.Range("E28").Value = "  -6.18%  "
$ws.Range("D29").Value = "'0.992"
$ws.Range("E29").Value = "  -1.04%  "
$ws.Range("E30").Value = "  -4.27%  "
$ws.Range("E31").Value = "  -4.40%  "
$ws.Range("E32").Value = "  -4.31%  "
$ws.Range("D33").Value = "'0.131"
$ws.Range("E33").Value = "  -5.36%  "
$ws.Range("D34").Value = "'0.999"
$ws.Range("E34").Value = "  +0.03%  "
$ws.Range("D35").Value = "'1.48"
$ws.Range("E35").Value = "  -5.40%  "
$ws.Range("D36").Value = "'4.64"
$ws.Range("E36").Value = "  -3.79%  "
$ws.Range("D37").Value = "'151.49"
$ws.Range("E37").Value = "  -1.72%  "
$ws.Range("D38").Value = "'0.363"
$ws.Range("E38").Value = "  -3.05%  "
$ws.Range("D39").Value = "'18.44"
$ws.Range("E40").Value = "  -5.92%  "
$ws.Range("D41").Value = "'2.68"
$ws.Range("E41").Value = "  -3.61%  "
$ws.Range("D44").Value = "0.0₆0307"
$ws.Range("E44").Value = "  +1.54%  "
$ws.Range("D45").Value = "'151.95"
$ws.Range("E45").Value = "  +2.45%  "
$ws.Range("D46").Value = "'15.30"
$ws.Range("E46").Value = "  +0.93%  "
$ws.Range("D47").Value = "'3.53"
$ws.Range("E47").Value = "  -3.01%  "
$ws.Range("E48").Value = "  -1.31%  "
$ws.Range("D49").Value = "'19.96"
$ws.Range("E49").Value = "  -4.92%  "
$ws.Range("E50").Value = "  -3.60%  "
$ws.Range("D51").Value = "'0.0900"
$ws.Range("E51").Value = "  -2.48%  "
